$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; F=1; G=1.467507333333333; H=4.402521999999999; I=0.1890754490804; J=0.1890754490804; K=3; L=1; M=0.2520896666666667; N=0.7562690000000001; O=0.03491140780587004; P=0.03491140780587004; Q=0.3699434344908889; R=3.329490910418; S=0.006600890108923859; T=0.006600890108923859 }
    3 = @{ E=3; F=1; G=1.467507333333333; H=4.402521999999999; I=0.1890754490804; J=0.1890754490804; K=3; L=1; M=1.312792666666667; N=3.938378; O=0.1818061039810792; P=0.1818061039810792; Q=1.926532865479555; R=17.338795789316; S=0.03437507075578045; T=0.03437507075578045 }
    4 = @{ E=3; F=1; G=1.467507333333333; H=4.402521999999999; I=0.1890754490804; J=0.1890754490804; K=3; L=1; M=5.655957; N=16.967871; O=0.7832824882130508; P=0.7832824882130508; Q=8.300158374517999; R=74.70142537066198; S=0.1480994882156957; T=0.1480994882156957 }
    5 = @{ E=3; F=1; G=6.293983333333333; H=18.88195; I=0.8109245509196; J=0.8109245509195999; K=3; L=1; M=0.2520896666666667; N=0.7562690000000001; O=0.03491140780587004; P=0.03491140780587004; Q=1.586648160505556; R=14.27983344455; S=0.02831051769694618; T=0.02831051769694617 }
    6 = @{ E=3; F=1; G=6.293983333333333; H=18.88195; I=0.8109245509196; J=0.8109245509195999; K=3; L=1; M=1.312792666666667; N=3.938378; O=0.1818061039810792; P=0.1818061039810792; Q=8.262695164122222; R=74.3642564771; S=0.1474310332252988; T=0.1474310332252987 }
    7 = @{ E=3; F=1; G=6.293983333333333; H=18.88195; I=0.8109245509196; J=0.8109245509195999; K=3; L=1; M=5.655957; N=16.967871; O=0.7832824882130508; P=0.7832824882130508; Q=35.59849909205; R=320.3864918284499; S=0.6351829999973552; T=0.6351829999973551 }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
